# Apply the "Metadata" sheet updates that correspond to a regenerated
# FHIR ValueSet spreadsheet export (Version bump, refreshed Date,
# a real Publisher, a new Jurisdiction row, and the removal of the
# duplicated "Contact" rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: refreshed export timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value
$ws.Range("B9").Value = "Alvearie Team"

# First "Contact" row becomes "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Second "Contact" row becomes "Description" (keeps the description text
# that used to live two rows further down)
$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "The value set that indicates the specific types of costs (admin fees, premiums, etc.) that may be tied to a member based on their Insurance plan cost category"

# Remaining rows shift up by one: Purpose, Copyright, Immutable
$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = ""
$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = ""
$ws.Range("A14").Value = "Immutable"
$ws.Range("B14").Value = "BooleanType[null]"

# The sheet now has one fewer row (A1:B14 instead of A1:B15); drop the
# trailing row entirely rather than leaving it blank.
$ws.Rows.Item(15).Delete()
